$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Nguon goc" (Source/Origin) header in column K, matching the
# formatting of the other single-line header cells (e.g. A1).
$ws.Range("K1").Value = "Nguồn gốc"

$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Match the new column's width to the rest of the header row style.
$ws.Range("K1").ColumnWidth = 18.43

# Move the active selection, as captured when the file was last saved.
$null = $ws.Range("I1").Select()
